$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new column before the current "Wellness" column (G) so the
# existing "State" column (F) is followed by a new "State Abbreviation"
# column, and the old "Wellness" column shifts right to H.
$ws.Columns("G").Insert()

# Rename the "State" header to "State Name" (values unchanged).
$ws.Range("F1").Value = "State Name"

# New "State Abbreviation" column.
$ws.Range("G1").Value = "State Abbreviation"
$ws.Range("G2").Value = "TX"
$ws.Range("G3").Value = "TX"

# "Wellness" column values change from "Ill" to "Well".
$ws.Range("H2").Value = "Well"
$ws.Range("H3").Value = "Well"

# Reset the selection to the top-left cell (clears the stale A2:G3
# selection rectangle left over from before the edit).
$ws.Range("A1").Select() | Out-Null
